$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1833.3334
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 1833.3334
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 1833.3334
$ws.Range("M40").Value = $null
$ws.Range("N40").Value = -2183.3334
$ws.Range("H62").Value = 5284.1665
$ws.Range("J62").Value = 6268.3335
$ws.Range("L62").Value = 6268.3335
$ws.Range("N62").Value = -7516.3335
$ws.Range("H65").Value = 5284.1665
$ws.Range("J65").Value = 6268.3335
$ws.Range("L65").Value = 31341.6675
$ws.Range("N65").Value = -37581.6675
$ws.Range("H76").Value = 3452.5
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 3452.5
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 3452.5
$ws.Range("M76").Value = $null
$ws.Range("N76").Value = -4082.5
$ws.Range("H79").Value = 3452.5
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 3452.5
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 3452.5
$ws.Range("M79").Value = $null
$ws.Range("N79").Value = -5636.5
$ws.Range("H106").Value = 4000
$ws.Range("I106").Value = 3000
$ws.Range("J106").Value = 4285.7144
$ws.Range("K106").Value = 3000
$ws.Range("L106").Value = 4285.7144
$ws.Range("M106").Value = -2369
$ws.Range("N106").Value = -5547.7144
$ws.Range("H112").Value = 1417.09
$ws.Range("I112").Value = 688.1667
$ws.Range("J112").Value = 1463.6171
$ws.Range("K112").Value = 2064.5001
$ws.Range("L112").Value = 4390.8513
$ws.Range("M112").Value = -956.5001000000002
$ws.Range("N112").Value = -6606.8513
$ws.Range("H115").Value = 1443.1875
$ws.Range("I115").Value = 519.1
$ws.Range("J115").Value = 2983.3333
$ws.Range("K115").Value = 1557.3
$ws.Range("L115").Value = 8949.999899999999
$ws.Range("M115").Value = 9.699999999999818
$ws.Range("N115").Value = -12083.9999
$ws.Range("H137").Value = 3230856.5
$ws.Range("I137").Value = 5269189
$ws.Range("J137").Value = 3497
$ws.Range("K137").Value = 15807567
$ws.Range("L137").Value = 10491
$ws.Range("M137").Value = -15805017
$ws.Range("N137").Value = -15591
$ws.Range("H138").Value = 7131.912
$ws.Range("I138").Value = 3350.7144
$ws.Range("K138").Value = 10052.1432
$ws.Range("M138").Value = -4912.143199999999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 34
$ws.Range("I37").Value = 34
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 34
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = 239
$ws.Range("N37").Value = $null
$ws.Range("H61").Value = 3742.44
$ws.Range("I61").Value = 3056.1
$ws.Range("J61").Value = 4200
$ws.Range("K61").Value = 3056.1
$ws.Range("L61").Value = 4200
$ws.Range("M61").Value = -2844.1
$ws.Range("N61").Value = -4624
$ws.Range("H132").Value = 2387.383
$ws.Range("I132").Value = 1676.2858
$ws.Range("J132").Value = 4461.4165
$ws.Range("K132").Value = 5028.857400000001
$ws.Range("L132").Value = 13384.2495
$ws.Range("M132").Value = -2498.857400000001
$ws.Range("N132").Value = -18444.2495
$ws.Range("H136").Value = 3742.44
$ws.Range("I136").Value = 3056.1
$ws.Range("J136").Value = 4200
$ws.Range("K136").Value = 9168.299999999999
$ws.Range("L136").Value = 12600
$ws.Range("M136").Value = -6618.299999999999
$ws.Range("N136").Value = -17700

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 23000
$ws.Range("J35").Value = 23000
$ws.Range("L35").Value = 23000
$ws.Range("N35").Value = -23620
$ws.Range("H86").Value = 59286.668
$ws.Range("I86").Value = 13000
$ws.Range("J86").Value = 65072.5
$ws.Range("K86").Value = 13000
$ws.Range("L86").Value = 65072.5
$ws.Range("M86").Value = -11877
$ws.Range("N86").Value = -67318.5
$ws.Range("H89").Value = 59286.668
$ws.Range("I89").Value = 13000
$ws.Range("J89").Value = 65072.5
$ws.Range("K89").Value = 65000
$ws.Range("L89").Value = 325362.5
$ws.Range("M89").Value = -59384
$ws.Range("N89").Value = -336594.5
$ws.Range("H127").Value = 30000
$ws.Range("J127").Value = 30000
$ws.Range("L127").Value = 30000
$ws.Range("N127").Value = -39920

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2860609.2
$ws.Range("I31").Value = 3450484.5
$ws.Range("J31").Value = 9546.666999999999
$ws.Range("K31").Value = 3450484.5
$ws.Range("L31").Value = 9546.666999999999
$ws.Range("M31").Value = -3450189.5
$ws.Range("N31").Value = -10136.667
$ws.Range("H34").Value = 2860609.2
$ws.Range("I34").Value = 3450484.5
$ws.Range("J34").Value = 9546.666999999999
$ws.Range("K34").Value = 3450484.5
$ws.Range("L34").Value = 9546.666999999999
$ws.Range("M34").Value = -3450282.5
$ws.Range("N34").Value = -9950.666999999999
$ws.Range("H141").Value = 31547.223
$ws.Range("J141").Value = 32020
$ws.Range("L141").Value = 32020
$ws.Range("N141").Value = -42380

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 5402.409
$ws.Range("J34").Value = 8408.071
$ws.Range("L34").Value = 25224.213
$ws.Range("N34").Value = -25392.213
$ws.Range("H39").Value = 2100
$ws.Range("H131").Value = 2217.1155
$ws.Range("I131").Value = 7600
$ws.Range("J131").Value = 1515
$ws.Range("K131").Value = 22800
$ws.Range("L131").Value = 4545
$ws.Range("M131").Value = -17760
$ws.Range("N131").Value = -14625

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4482.8423
$ws.Range("I70").Value = 4459.5386
$ws.Range("J70").Value = 4533.3335
$ws.Range("K70").Value = 4459.5386
$ws.Range("L70").Value = 4533.3335
$ws.Range("M70").Value = -4189.5386
$ws.Range("N70").Value = -5073.3335
$ws.Range("H73").Value = 4482.8423
$ws.Range("I73").Value = 4459.5386
$ws.Range("J73").Value = 4533.3335
$ws.Range("K73").Value = 4459.5386
$ws.Range("L73").Value = 4533.3335
$ws.Range("M73").Value = -3523.5386
$ws.Range("N73").Value = -6405.3335

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 401000
$ws.Range("J2").Value = 5000
$ws.Range("L2").Value = 5000
$ws.Range("N2").Value = -5224
$ws.Range("H14").Value = 11611.889
$ws.Range("J14").Value = 11706.471
$ws.Range("L14").Value = 11706.471
$ws.Range("N14").Value = -12050.471
$ws.Range("H20").Value = 68006
$ws.Range("J20").Value = 68006
$ws.Range("L20").Value = 68006
$ws.Range("N20").Value = -68458
$ws.Range("H21").Value = 60000
$ws.Range("J21").Value = 60000
$ws.Range("L21").Value = 60000
$ws.Range("N21").Value = -60348
$ws.Range("H22").Value = 2400
$ws.Range("J22").Value = 2780
$ws.Range("L22").Value = 2780
$ws.Range("N22").Value = -3370
$ws.Range("H27").Value = 2400
$ws.Range("J27").Value = 2780
$ws.Range("L27").Value = 2780
$ws.Range("N27").Value = -2994
$ws.Range("H46").Value = 1158.4717
$ws.Range("I46").Value = 983.3333
$ws.Range("J46").Value = 2839.8
$ws.Range("K46").Value = 983.3333
$ws.Range("L46").Value = 2839.8
$ws.Range("M46").Value = -795.3333
$ws.Range("N46").Value = -3215.8
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").Value = $null
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").Value = $null
$ws.Range("H54").Value = 13695
$ws.Range("J54").Value = 13695
$ws.Range("L54").Value = 13695
$ws.Range("N54").Value = -14983
$ws.Range("H55").Value = 795.619
$ws.Range("I55").Value = 235.41667
$ws.Range("J55").Value = 1542.5555
$ws.Range("K55").Value = 235.41667
$ws.Range("L55").Value = 1542.5555
$ws.Range("M55").Value = -62.41667000000001
$ws.Range("N55").Value = -1888.5555
$ws.Range("H100").Value = 3650
$ws.Range("I100").Value = 1733.3334
$ws.Range("J100").Value = 4225
$ws.Range("K100").Value = 1733.3334
$ws.Range("L100").Value = 4225
$ws.Range("M100").Value = -1192.3334
$ws.Range("N100").Value = -5307

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1613.8
$ws.Range("I107").Value = 325.4
$ws.Range("K107").Value = 976.1999999999999
$ws.Range("M107").Value = 943.8000000000001
$ws.Range("H138").Value = 29732.25
$ws.Range("J138").Value = 29732.25
$ws.Range("L138").Value = 29732.25
$ws.Range("N138").Value = -40012.25
